$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 97
$ws.Cells.Item($newRow, 1).Value = 41
$ws.Cells.Item($newRow, 2).Value = 15
$ws.Cells.Item($newRow, 3).Value = 1.5
$ws.Cells.Item($newRow, 4).Value = 50
$ws.Cells.Item($newRow, 5).Value = 84.84999999999999
$ws.Cells.Item($newRow, 6).Value = 10201
